$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 40866162
$ws.Range("D2").Value = 34910128
$ws.Range("I2").Value = 25.2844034213071
$ws.Range("J2").Value = 25.0571423108366
$ws.Range("O2").Value = 25.2992826333767
$ws.Range("P2").Value = 25.1168400072513
$ws.Range("Q2").Value = 1.91835833608368
$ws.Range("R2").Value = 0.152298309114685
$ws.Range("S2").Value = 0.761491545573424
$ws.Range("T2").Value = 1.13180608072854
$ws.Range("U2").Value = 0.883543583151899
$ws.Range("V2").Value = -0.122879738982177
$ws.Range("W2").Value = 0.487764991232885
$ws.Range("X2").Value = 3
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = "#N/A"
$ws.Range("C3").Value = 200506080
$ws.Range("D3").Value = 228724928
$ws.Range("I3").Value = 27.5790707436736
$ws.Range("J3").Value = 27.769038368146
$ws.Range("O3").Value = 27.7551588118348
$ws.Range("P3").Value = 27.7095536145082
$ws.Range("Q3").Value = 0.628467294661039
$ws.Range("R3").Value = 0.575887683391777
$ws.Range("S3").Value = 0.938111873128734
$ws.Range("T3").Value = 1.03052639735946
$ws.Range("U3").Value = 0.970377859861055
$ws.Range("V3").Value = -0.190521784441083
$ws.Range("W3").Value = 0.281732179094314
$ws.Range("X3").Value = 3
$ws.Range("Y3").Value = 3
$ws.Range("Z3").Value = "#N/A"
$ws.Range("C4").Value = 26963010
$ws.Range("D4").Value = 25142570
$ws.Range("I4").Value = 24.6844782241466
$ws.Range("J4").Value = 24.5836287895766
$ws.Range("O4").Value = 25.0746667019424
$ws.Range("P4").Value = 24.5976798379387
$ws.Range("Q4").Value = 0.882504881863556
$ws.Range("R4").Value = 0.469411733074466
$ws.Range("S4").Value = 0.938111873128734
$ws.Range("T4").Value = 1.61385302131594
$ws.Range("U4").Value = 0.619635113478052
$ws.Range("V4").Value = -1.81599888633545
$ws.Range("W4").Value = 2.76997261434286
$ws.Range("X4").Value = 3
$ws.Range("Y4").Value = 3
$ws.Range("Z4").Value = "#N/A"
$ws.Range("H5").Value = 51412512
$ws.Range("N5").Value = 25.6156161676712
$ws.Range("C7").Value = 21140340
$ws.Range("D7").Value = 23729240
$ws.Range("I7").Value = 24.333495243977
$ws.Range("J7").Value = 24.5001625588526
$ws.Range("O7").Value = 24.4045916014765
$ws.Range("P7").Value = 24.4172683211009
$ws.Range("Q7").Value = -0.188190056380374
$ws.Range("R7").Value = 0.859892871192765
$ws.Range("S7").Value = 0.938111873128734
$ws.Range("T7").Value = 0.991205494871482
$ws.Range("U7").Value = 1.00887253468027
$ws.Range("V7").Value = -0.199732228740137
$ws.Range("W7").Value = 0.174378789491213
$ws.Range("X7").Value = 3
$ws.Range("Y7").Value = 3
$ws.Range("Z7").Value = "#N/A"
$ws.Range("C8").Value = 29319660
$ws.Range("D8").Value = 36068288
$ws.Range("I8").Value = 24.8053650378348
$ws.Range("J8").Value = 25.1042276106936
$ws.Range("O8").Value = 25.2308048413939
$ws.Range("P8").Value = 25.0359941775666
$ws.Range("Q8").Value = 1.59667575311067
$ws.Range("R8").Value = 0.238062091527051
$ws.Range("S8").Value = 0.793540305090171
$ws.Range("T8").Value = 1.13751585877202
$ws.Range("U8").Value = 0.879108622783973
$ws.Range("V8").Value = -0.279532260186281
$ws.Range("W8").Value = 0.669153587840989
$ws.Range("X8").Value = 3
$ws.Range("Y8").Value = 3
$ws.Range("Z8").Value = "#N/A"
$ws.Range("C9").Value = 34100360
$ws.Range("D9").Value = 37990528
$ws.Range("I9").Value = 25.0232836341787
$ws.Range("J9").Value = 25.1791364272197
$ws.Range("O9").Value = 25.0624917054832
$ws.Range("P9").Value = 25.054395391038
$ws.Range("Q9").Value = 0.082707818416188
$ws.Range("R9").Value = 0.938111873128734
$ws.Range("S9").Value = 0.938111873128734
$ws.Range("T9").Value = 1.00620594137914
$ws.Range("U9").Value = 0.993832334789599
$ws.Range("V9").Value = -0.26524486500591
$ws.Range("W9").Value = 0.281437493896425
$ws.Range("X9").Value = 3
$ws.Range("Y9").Value = 3
$ws.Range("Z9").Value = "#N/A"
$ws.Range("C10").Value = 2663829782
$ws.Range("D10").Value = 2590670618
$ws.Range("I10").Value = 31.3108547515192
$ws.Range("J10").Value = 31.2706784546122
$ws.Range("O10").Value = 31.2119217219633
$ws.Range("P10").Value = 31.3052217494229
$ws.Range("Q10").Value = -4.40168491352323
$ws.Range("R10").Value = 0.0199533891722438
$ws.Range("S10").Value = 0.199533891722438
$ws.Range("T10").Value = 0.937269329894508
$ws.Range("U10").Value = 1.06692918257824
$ws.Range("V10").Value = -0.159248621463293
$ws.Range("W10").Value = -0.0273514334559108
$ws.Range("X10").Value = 3
$ws.Range("Y10").Value = 3
$ws.Range("Z10").Value = "#N/A"
$ws.Range("C11").Value = "#N/A"
$ws.Range("D11").Value = 1495972
$ws.Range("I11").Value = "#N/A"
$ws.Range("J11").Value = 20.5126517419834
$ws.Range("C12").Value = 36614408
$ws.Range("D12").Value = 52227550
$ws.Range("I12").Value = 25.1259081340278
$ws.Range("J12").Value = 25.6383076925291
$ws.Range("O12").Value = 25.526739600459
$ws.Range("P12").Value = 25.4670813316822
$ws.Range("Q12").Value = 0.347529614872081
$ws.Range("R12").Value = 0.760632790466221
$ws.Range("S12").Value = 0.938111873128734
$ws.Range("T12").Value = 1.02854192623387
$ws.Range("U12").Value = 0.972250109105049
$ws.Range("V12").Value = -0.661741153627516
$ws.Range("W12").Value = 0.781057691181139
$ws.Range("X12").Value = 3
$ws.Range("Y12").Value = 3
$ws.Range("Z12").Value = "#N/A"
$ws.Range("C14").Value = 49490460
$ws.Range("D14").Value = 49514448
$ws.Range("I14").Value = 25.5606471159261
$ws.Range("J14").Value = 25.5613462200416
$ws.Range("O14").Value = 25.570656813745
$ws.Range("P14").Value = 25.5658457330176
$ws.Range("Q14").Value = 0.120093682070113
$ws.Range("R14").Value = 0.915103230244877
$ws.Range("S14").Value = 0.938111873128734
$ws.Range("T14").Value = 1.00408756419309
$ws.Range("U14").Value = 0.99592907597021
$ws.Range("V14").Value = -0.162867450169104
$ws.Range("W14").Value = 0.172489611623869
$ws.Range("X14").Value = 3
$ws.Range("Y14").Value = 3
$ws.Range("Z14").Value = "#N/A"
$ws.Range("C15").Value = "#N/A"
$ws.Range("D15").Value = 15940000
$ws.Range("I15").Value = "#N/A"
$ws.Range("J15").Value = 23.9261482935399
$ws.Range("C16").Value = 64759632
$ws.Range("D16").Value = 53367728
$ws.Range("I16").Value = 25.9485914516852
$ws.Range("J16").Value = 25.6694642576456
$ws.Range("O16").Value = 25.8567947696956
$ws.Range("P16").Value = 25.8214960721334
$ws.Range("Q16").Value = 0.408847988878684
$ws.Range("R16").Value = 0.715351772062022
$ws.Range("S16").Value = 0.938111873128734
$ws.Range("T16").Value = 1.02192716928869
$ws.Range("U16").Value = 0.978543315074051
$ws.Range("V16").Value = -0.275304872879686
$ws.Range("W16").Value = 0.34590226800396
$ws.Range("X16").Value = 3
$ws.Range("Y16").Value = 3
$ws.Range("Z16").Value = "#N/A"
